$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("B2").Value = 17
$ws.Range("B5").Value = 6
